$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two "CD55 / Gene" rows (original rows 6 and 7)
$ws.Rows("6:7").Delete()

# Relabel the "Keyword" column values from "Metastasis" to "Cancer"
# (now rows 2-5 after the deletion above)
$ws.Range("D2:D5").Value = "Cancer"

# Move the active selection to F12, matching the saved selection in the workbook
$ws.Range("F12").Select()
